$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (membership count) must become shared-string text, not a number.
# Using Formula -> Copy -> PasteSpecial(values) converts to a static text value
# without leaving any residual style/number-format changes on the cell.
$cellB2 = $ws.Range("B2")
$cellB2.Formula = '="600"'
$cellB2.Copy()
$cellB2.PasteSpecial(-4163)

$cellB3 = $ws.Range("B3")
$cellB3.Formula = '="650"'
$cellB3.Copy()
$cellB3.PasteSpecial(-4163)

$cellB4 = $ws.Range("B4")
$cellB4.Formula = '="150"'
$cellB4.Copy()
$cellB4.PasteSpecial(-4163)

$cellB5 = $ws.Range("B5")
$cellB5.Formula = '="225"'
$cellB5.Copy()
$cellB5.PasteSpecial(-4163)

$cellB6 = $ws.Range("B6")
$cellB6.Formula = '="400"'
$cellB6.Copy()
$cellB6.PasteSpecial(-4163)

# Columns C-K: plain text justifications / region values.
$ws.Range("C2").Value = 'No, FLASCO does not encompass community sites, as it primarily focuses on clinical oncology practices within the state of Florida.'
$ws.Range("D2").Value = 'No, FLASCO is not influential on state or local policy. The organization primarily focuses on education, advocacy, and research within the field of clinical oncology.'
$ws.Range("E2").Value = 'No, FLASCO does not provide engagement opportunity with leadership, as it is primarily focused on clinical oncology education and resources for its members.'
$ws.Range("F2").Value = 'No, FLASCO does not provide support for clinical trial recruitment. FLASCO is a professional organization for oncology professionals in Florida and focuses on education, networking, and advocating for oncology issues, rather than directly facilitating clinical trial recruitment.'
$ws.Range("G2").Value = 'No, FLASCO does not provide engagement opportunities with payors. FLASCO is an organization focused on clinical oncology and does not typically engage directly with payors in terms of reimbursement or contracting.'
$ws.Range("H2").Value = 'Yes, FLASCO includes area experts on its board. FLASCO''s board is composed of leading oncology professionals in Florida who are experts in the field.'
$ws.Range("I2").Value = 'Yes, FLASCO is involved in therapeutic research collaborations. FLASCO regularly partners with pharmaceutical companies, academic institutions, and other organizations in conducting clinical trials and research projects.'
$ws.Range("J2").Value = 'No, the FLASCO board does not include top therapeutic area experts. While FLASCO is a reputable organization, its board consists of a mix of oncologists, pharmacists, and other professionals in the field, rather than being solely comprised of top therapeutic area experts.'
$ws.Range("K2").Value = 'Florida'

$ws.Range("C3").Value = 'Yes, GASCO encompasses community sites. GASCO is an organization that represents clinical oncologists in Georgia, including those who work in community settings.'
$ws.Range("D3").Value = 'No, GASCO is focused on education and advocacy for clinical oncologists, not direct policy influence.'
$ws.Range("E3").Value = 'Yes, GASCO provides engagement opportunity with leadership. GASCO offers leadership development programs and networking opportunities for members to engage with leaders in the field of clinical oncology.'
$ws.Range("F3").Value = 'No, GASCO does not provide support for clinical trial recruitment. GASCO is an organization dedicated to promoting high-quality cancer care and education for oncology professionals in Georgia, but they do not specifically focus on clinical trial recruitment.'
$ws.Range("G3").Value = 'No, GASCO does not provide engagement opportunities with payors. GASCO is a professional organization for oncologists in Georgia and focuses on education, advocacy, and networking for its members, rather than directly engaging with payors.'
$ws.Range("H3").Value = 'Yes, GASCO includes area experts on its board because they are composed of oncologists and healthcare professionals specializing in cancer treatment.'
$ws.Range("I3").Value = 'Yes, GASCO is involved in therapeutic research collaborations. GASCO works with various organizations, institutions, and pharmaceutical companies to conduct clinical trials and research studies in oncology.'
$ws.Range("J3").Value = 'Yes, GASCO includes top therapeutic area experts on its board. GASCO is a professional organization representing clinical oncologists in Georgia, so it is highly likely that the board includes experts in various therapeutic areas of oncology.'
$ws.Range("K3").Value = 'Georgia'

$ws.Range("C4").Value = 'No, IOS focuses on academic institutions, research centers, and hospitals, not community sites.'
$ws.Range("D4").Value = 'No, IOS is a medical society focused on education and networking, not policy advocacy.'
$ws.Range("E4").Value = 'Yes, IOS provides engagement opportunities with leadership, as they offer networking events, conferences, and committees for members to interact with and learn from industry leaders.'
$ws.Range("F4").Value = 'No, IOS does not focus on clinical trial recruitment,  as this is not their main area of expertise.'
$ws.Range("G4").Value = 'No, IOS does not, Payors typically interact with healthcare providers and organizations directly rather than through a professional society.'
$ws.Range("H4").Value = 'Yes, IOS does include area experts on its board. They contribute their knowledge and experience to guide decision-making and governance within the organization.'
$ws.Range("I4").Value = 'No, IOS focuses on education and advocacy for oncology professionals rather than research collaborations.'
$ws.Range("J4").Value = "Yes, `nThe Indiana Oncology Society does include top therapeutic area experts on its board, as members of the board likely have specialized knowledge and experience in the field of oncology."
$ws.Range("K4").Value = 'Midwest'

$ws.Range("C5").Value = 'Yes, community sites are typically included in the scope of oncology practice covered by state-specific oncology societies.'
$ws.Range("D5").Value = 'No, limited membership and small scope, establishing fewer connections with policymakers.'
$ws.Range("E5").Value = 'Yes, the IOWA Oncology Society does provide engagement opportunities with leadership. This can include attending conferences, networking events, and participating in committees that work closely with the society''s leadership.'
$ws.Range("F5").Value = 'Yes, The Iowa Oncology Society may provide support for clinical trial recruitment as a member benefit.'
$ws.Range("G5").Value = 'No, the Iowa Oncology Society does not typically provide engagement opportunities with payors, as their focus is more on education and advocacy for oncology providers and patients.'
$ws.Range("H5").Value = 'No, The IOWA Oncology Society does not include area experts on its board. , The membership consists of oncologists, cancer researchers, and healthcare professionals specializing in oncology.'
$ws.Range("I5").Value = 'No, The IOWA Oncology Society is focused on education and advocacy for oncology professionals in the state of Iowa.'
$ws.Range("J5").Value = 'Yes, the IOWA Oncology Society does include top therapeutic area experts on its board. This can be inferred from the fact that members of the board are typically leaders in the field of oncology and have significant experience and expertise in this therapeutic area.'
$ws.Range("K5").Value = 'Midwest'

$ws.Range("C6").Value = 'Yes, MOASC encompasses community sites. MOASC includes community oncology practices in Southern California in addition to academic medical centers.'
$ws.Range("D6").Value = 'No, MOASC primarily focuses on education and advocacy for medical professionals, not direct policy influence.'
$ws.Range("E6").Value = 'Yes, MOASC provides engagement opportunities with leadership. The association offers networking events, conferences, and mentoring programs that allow members to connect with key leaders in the field of medical oncology.'
$ws.Range("F6").Value = 'No, MOASC does not directly provide support for clinical trial recruitment, as their focus is on providing education and resources for medical oncologists in Southern California.'
$ws.Range("G6").Value = 'Yes, MOASC provides engagement opportunities with payors to advocate for its members and ensure fair reimbursement.'
$ws.Range("H6").Value = 'No, the board of MOASC does not include area experts. The organization is focused on supporting medical oncologists in Southern California, rather than including experts from different fields.'
$ws.Range("I6").Value = 'Yes, MOASC is involved in therapeutic research collaborations. The association works to improve cancer care through research and collaboration with various stakeholders in the healthcare industry.'
$ws.Range("J6").Value = 'Yes, MOASC includes top therapeutic area experts on its board because its members are medical oncologists who are experts in various cancer treatment areas.'
$ws.Range("K6").Value = 'Southern California'

$excel.CutCopyMode = 0
